$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Russian header in B1: "1.3.1.1.f" -> "1.3.1.1f" (drop the stray dot)
$ws.Range("B1").Value = "1.3.1.1f Доля лиц, получающих пенсии и пособия по инвалидности к общей численности населения"

# Row height tweaks that came with the new column
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 29.25

# Add a new year column (T) mirroring column S's formatting, then fill in the 2023 data
$ws.Range("S2:S5").Copy($ws.Range("T2:T5"))
$ws.Range("T2").Value = ""
$ws.Range("T3").Value = 2023
$ws.Range("T4").Value = 217222
$ws.Range("T5").Value = 2.9794303052841493
